$d = $word.ActiveDocument
$wmlns = ' xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1) Remove the old "_GoBack" bookmark that sits alone in the page-break
#    paragraph right before the "Outil et source" heading.
#    We locate that paragraph via the page-break that immediately precedes
#    the heading text "Outil ", then rewrite its (tiny) range twice:
#    first collapse it down to nothing (this strips the bookmark + run),
#    then restore the page-break run + paragraph properties, this time
#    without re-creating the bookmark.
# ---------------------------------------------------------------------------

$rng = $d.Content
$found = $rng.Find.Execute("Outil ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$headingStart = $rng.Start

# the paragraph holding the lone page break ends exactly one character
# (its own paragraph mark) before the heading paragraph starts
$pbRange = $d.Range($headingStart - 2, $headingStart - 1)
$null = $pbRange.InsertXML('<w:r' + $wmlns + '><w:t></w:t></w:r>')

$pbRange2 = $d.Range($headingStart - 2, $headingStart - 1)
$null = $pbRange2.InsertXML('<w:p' + $wmlns + '><w:pPr><w:spacing w:after="160"/><w:jc w:val="left"/></w:pPr><w:r><w:br w:type="page"/></w:r></w:p>')

# ---------------------------------------------------------------------------
# 2) Append the new "PROBA / MATRICE / GRAPH" reference section at the very
#    end of the document body (before the final section properties).
# ---------------------------------------------------------------------------

$newContent = @'
<w:p/><w:p/><w:p><w:pPr><w:spacing w:after="160"/><w:jc w:val="left"/></w:pPr><w:r><w:br w:type="page"/></w:r></w:p><w:p><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">PROBA </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>GenererAgeMoy</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> l 80</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Fleuve 149 &#8211; 202</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>zoneUrbaine</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> 215-239</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:t>composition</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> de grille 355-407</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>genererZoneUrbaine</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> 418 455</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>lancerVoyage</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/></w:p><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>soumettreVirus</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/></w:p><w:p/><w:p><w:r><w:t>MATRICE</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>genererDeplacment</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> 458</w:t></w:r></w:p><w:p><w:bookmarkStart w:id="401" w:name="_GoBack"/><w:bookmarkEnd w:id="401"/></w:p><w:p><w:r><w:t>GRAPH</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>algoPrim</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> 1023</w:t></w:r></w:p>
'@
$newContent = $newContent.Replace('<w:p>', '<w:p' + $wmlns + '>').Replace('<w:p/>', '<w:p' + $wmlns + '/>')

$endRng = $d.Content
$endRng.Collapse(0)
$null = $endRng.InsertXML($newContent)

Write-Output "done"
